# Reorder the "Salida" values (column B) on the "Solucion" worksheet so that
# the shared-string table ends up rebuilt in the order required by the diff.
# Column A ("Pedido") values stay untouched; only column B values change,
# which is exactly the visible effect of the shared-strings reshuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solucion")

$salidas = @(
    "S001", "S021", "S041", "S051", "S031", "S011", "S042", "S002", "S022", "S052",
    "S012", "S032", "S043", "S003", "S023", "S053", "S033", "S013", "S044", "S004",
    "S024", "S034", "S014", "S054", "S025", "S045", "S005", "S055", "S015", "S035",
    "S026", "S046", "S006", "S056", "S036", "S016", "S047", "S007", "S027", "S057",
    "S017", "S037", "S028", "S048", "S008", "S038", "S058", "S018", "S009", "S049",
    "S029", "S059", "S039", "S019", "S030", "S050", "S010", "S060", "S020", "S040"
)

for ($i = 0; $i -lt $salidas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $salidas[$i]
}
